$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# New test case row (Watchlist034) appended below the existing data (row 34 -> row 35).
# Copy the formatting of row 27 (style pattern 7/1/1/1/1 across A:E, matching what row 35
# needs) onto the new row before writing values, so no new cell-style entries are created.
$ws.Range("A27:E27").Copy()
$ws.Range("A35:E35").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row - Jira id / Description / TCID order matches the original
# shared-string allocation order captured in the target workbook.
$ws.Cells.Item(35, 2).Value = "OPQA-322"
$ws.Cells.Item(35, 3).Value = "Verify that user is able to comment on his watchlist items"
$ws.Cells.Item(35, 1).Value = "Watchlist034"
$ws.Cells.Item(35, 4).Value = "Y"
$ws.Cells.Item(35, 5).Value = ""

# Update the sheet view: scroll the window so row 31 / column C is the top-left visible
# cell, and move the active selection to C37.
$win = $excel.ActiveWindow
$win.ScrollRow = 31
$win.ScrollColumn = 3
$ws.Range("C37").Select()
